$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 375 (shifts existing rows 375..401 down to 376..402)
$ws.Rows(375).Insert()

# Populate the newly inserted row with the new "Ají" price record
$ws.Cells.Item(375, 1).Value = 9
$ws.Cells.Item(375, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(375, 3).Value = "Metropolitana"
$ws.Cells.Item(375, 4).Value = 45041
$ws.Cells.Item(375, 5).Value = 13
$ws.Cells.Item(375, 6).Value = 100112021
$ws.Cells.Item(375, 7).Value = "Ají"
$ws.Cells.Item(375, 8).Value = "Americana (o)"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 52
$ws.Cells.Item(375, 11).Value = 13000
$ws.Cells.Item(375, 12).Value = 15000
$ws.Cells.Item(375, 13).Value = 14000
$ws.Cells.Item(375, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(375, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(375, 16).Value = 560
$ws.Cells.Item(375, 17).Value = 25
$ws.Cells.Item(375, 18).Value = "Hortaliza"

# Match the date formatting used by the rest of column D
$ws.Range("D375").NumberFormat = $ws.Range("D376").NumberFormat
